$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.43%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'40.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'2.56%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.130"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.80%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07614"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.06%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.624"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.04%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'2.437"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.01%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8989"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'2.46%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1084"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'12.00%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1768"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.40%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09227"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'3.27%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04202"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-4.56%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1050"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.61%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001252"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.62%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005879"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.34%"
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'-0.13%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'4.251"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.60%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D19").Value = "'6.573"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-6.07%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'1.86%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2681"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-14.54%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04092"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-1.76%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001224"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'2.41%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004087"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.58%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'6.64%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D38").Value = "'0.02371"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'1.22%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05175"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'0.54%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007781"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.79%"
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'-1.82%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'5.96%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.001952"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'0.13%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008562"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-0.42%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3073"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.82%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006925"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'6.11%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.01%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.03376"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'895.03%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.004201"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-40.00%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.01%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.01%"
$ws.Range("E51").Style = "Normal"
